$wb = $excel.ActiveWorkbook

# --- Sheet1: Schedule ---
$ws1 = $wb.Worksheets.Item("Schedule")

$ws1.Range("B4").Value = 46070.125
$ws1.Range("C4").Value = 4.5
$ws1.Range("D4").Value = 17.01
$ws1.Range("E4").Value = 502.9321635
$ws1.Range("F4").Value = 29.56685264550265
$ws1.Range("A5").Value = 46070.29166666666
$ws1.Range("C5").Value = 9
$ws1.Range("D5").Value = 34.02
$ws1.Range("E5").Value = 47.57589524999998
$ws1.Range("F5").Value = 1.398468408289241

# --- Sheet2: Detailed ---
$ws2 = $wb.Worksheets.Item("Detailed")

$ws2.Range("B38").Value = 73.01385000000001
$ws2.Range("B39").Value = 74.76455
$ws2.Range("B40").Value = 90.56785000000001
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 84.79000000000001
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 84.79000000000001
$ws2.Range("C42").Value = "historical"
$ws2.Range("B43").Value = 77.93199
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 72.8985
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 70.74718
$ws2.Range("C45").Value = "historical"
$ws2.Range("B46").Value = 62.65424
$ws2.Range("C46").Value = "historical"
$ws2.Range("B47").Value = 63.42016
$ws2.Range("C47").Value = "historical"
$ws2.Range("B48").Value = 57.06012
$ws2.Range("C48").Value = "historical"
$ws2.Range("C49").Value = "historical"
$ws2.Range("B52").Value = 53.46758
$ws2.Range("B53").Value = 56.98
$ws2.Range("B54").Value = 56.98
$ws2.Range("B55").Value = 56.98
$ws2.Range("B56").Value = 56.98
$ws2.Range("E56").Value = "OFF"
$ws2.Range("B57").Value = 56.98
$ws2.Range("B58").Value = 57.06007
$ws2.Range("B59").Value = 57.31
$ws2.Range("B60").Value = 64.83502
$ws2.Range("B62").Value = 75.83856
$ws2.Range("B63").Value = 69.92309
$ws2.Range("B64").Value = 56.98
$ws2.Range("E64").Value = "ON"
$ws2.Range("B65").Value = 28.67165
$ws2.Range("B66").Value = 14.57277
$ws2.Range("B67").Value = 0.51
$ws2.Range("B68").Value = 0.26957
$ws2.Range("B69").Value = -4.53645
$ws2.Range("B70").Value = -5.97422
$ws2.Range("B71").Value = -5.14767
$ws2.Range("B72").Value = -5.88559
$ws2.Range("B73").Value = -5.96707
$ws2.Range("B74").Value = -5.73682
$ws2.Range("B75").Value = -8.500590000000001
$ws2.Range("B76").Value = -8.426539999999999
$ws2.Range("B77").Value = -8.19598
$ws2.Range("B78").Value = -4.08444
$ws2.Range("B79").Value = -1.78797
$ws2.Range("B80").Value = 4.87821
$ws2.Range("B81").Value = 7.15693
$ws2.Range("B82").Value = 5.91348
$ws2.Range("B84").Value = 59.71376
$ws2.Range("B85").Value = 65.54300000000001
$ws2.Range("B86").Value = 78
$ws2.Range("B87").Value = 86.2732
$ws2.Range("B88").Value = 107.07418
$ws2.Range("B89").Value = 107.18745
$ws2.Range("B90").Value = 99.57449
$ws2.Range("B91").Value = 85.03959999999999
$ws2.Range("B92").Value = 95.16157
$ws2.Range("B93").Value = 84.79000000000001
$ws2.Range("B94").Value = 78
$ws2.Range("B95").Value = 73.2
$ws2.Range("B96").Value = 74.50967
$ws2.Range("B97").Value = 76.02867000000001

Write-Output "edits applied"
